$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 660
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H43").Value = 1006.9167
$ws.Range("I43").Value = 1091.875
$ws.Range("J43").Value = 837
$ws.Range("K43").Value = 1091.875
$ws.Range("L43").Value = 837
$ws.Range("M43").Value = -1022.875
$ws.Range("N43").Value = -975

$ws.Range("H96").Value = 7962.4
$ws.Range("J96").Value = 1878
$ws.Range("L96").Value = 5634
$ws.Range("N96").Value = -8380

$ws.Range("H97").Value = 999.6667
$ws.Range("J97").Value = 999.6667
$ws.Range("L97").Value = 2999.0001
$ws.Range("N97").Value = -3991.0001

$ws.Range("H100").Value = 999.5
$ws.Range("I100").Value = 999
$ws.Range("K100").Value = 999
$ws.Range("M100").Value = -458

$ws.Range("H113").Value = 2381.6924
$ws.Range("I113").Value = 1665.8889
$ws.Range("J113").Value = 3992.25
$ws.Range("K113").Value = 1665.8889
$ws.Range("L113").Value = 3992.25
$ws.Range("M113").Value = 1588.1111
$ws.Range("N113").Value = -10500.25

$ws.Range("H129").Value = 3185
$ws.Range("I129").Value = 2197
$ws.Range("K129").Value = 6591
$ws.Range("M129").Value = -1591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1858
$ws.Range("I74").Value = 1667.6666
$ws.Range("K74").Value = 1667.6666
$ws.Range("M74").Value = -793.6666

$ws.Range("H77").Value = 1858
$ws.Range("I77").Value = 1667.6666
$ws.Range("K77").Value = 8338.333000000001
$ws.Range("M77").Value = -3970.333000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 203.41667
$ws.Range("I22").Value = 204.75
$ws.Range("K22").Value = 204.75
$ws.Range("M22").Value = -31.75

$ws.Range("H86").Value = 1057
$ws.Range("I86").Value = 850.1429000000001
$ws.Range("K86").Value = 850.1429000000001
$ws.Range("M86").Value = 272.8570999999999

$ws.Range("H89").Value = 1057
$ws.Range("I89").Value = 850.1429000000001
$ws.Range("K89").Value = 4250.7145
$ws.Range("M89").Value = 1365.2855

$ws.Range("H94").Value = 2142.7856
$ws.Range("I94").Value = 2083.25
$ws.Range("K94").Value = 2083.25
$ws.Range("M94").Value = -1632.25

$ws.Range("H99").Value = 1989
$ws.Range("J99").Value = 1989
$ws.Range("L99").Value = 1989
$ws.Range("N99").Value = -4985

$ws.Range("H134").Value = 2190.125
$ws.Range("I134").Value = 2217.2856
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 6651.8568
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -4116.8568
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 935.4
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 959.6667
$ws.Range("K22").Value = 899
$ws.Range("L22").Value = 959.6667
$ws.Range("M22").Value = -549
$ws.Range("N22").Value = -1659.6667

$ws.Range("H86").Value = 8657.143
$ws.Range("I86").Value = 7319.2
$ws.Range("K86").Value = 7319.2
$ws.Range("M86").Value = -6196.2

$ws.Range("H89").Value = 8657.143
$ws.Range("I89").Value = 7319.2
$ws.Range("K89").Value = 36596
$ws.Range("M89").Value = -30980

$ws.Range("H132").Value = 6567.8184
$ws.Range("I132").Value = 6224.7
$ws.Range("K132").Value = 18674.1
$ws.Range("M132").Value = -16144.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 368547.66
$ws.Range("I2").Value = 366713.34
$ws.Range("K2").Value = 2200280.04
$ws.Range("M2").Value = -2200167.04

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

$ws.Range("H12").Value = 128.8
$ws.Range("I12").Value = 184.22223
$ws.Range("K12").Value = 552.66669
$ws.Range("M12").Value = -379.66669

$ws.Range("H22").Value = 1928.6923
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 2107.3
$ws.Range("K22").Value = 4000.0002
$ws.Range("L22").Value = 6321.900000000001
$ws.Range("M22").Value = -3831.0002
$ws.Range("N22").Value = -6659.900000000001

$ws.Range("H27").Value = 1928.6923
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 2107.3
$ws.Range("K27").Value = 4000.0002
$ws.Range("L27").Value = 6321.900000000001
$ws.Range("M27").Value = -3898.0002
$ws.Range("N27").Value = -6525.900000000001

$ws.Range("H32").Value = 549.5

$ws.Range("H34").Value = 1361
$ws.Range("J34").Value = 1361
$ws.Range("L34").Value = 4083
$ws.Range("N34").Value = -4251

$ws.Range("H38").Value = 3574.2
$ws.Range("I38").Value = 3574.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 10722.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -10375.6
$ws.Range("N38").ClearContents()

$ws.Range("H44").Value = 198.33333
$ws.Range("I44").Value = 225
$ws.Range("J44").Value = 190.71428
$ws.Range("K44").Value = 675
$ws.Range("L44").Value = 572.14284
$ws.Range("M44").Value = -277
$ws.Range("N44").Value = -1368.14284

$ws.Range("H46").Value = 1492.6666
$ws.Range("J46").Value = 1492.6666
$ws.Range("L46").Value = 4477.9998
$ws.Range("N46").Value = -4659.9998

$ws.Range("H51").Value = 2499
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H95").Value = 8873.666999999999
$ws.Range("J95").Value = 8873.666999999999
$ws.Range("L95").Value = 26621.001
$ws.Range("N95").Value = -30739.001

$ws.Range("H97").Value = 1896
$ws.Range("I97").Value = 1896
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5688
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -5192
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2175167.2
$ws.Range("J11").Value = 757500
$ws.Range("L11").Value = 757500
$ws.Range("N11").Value = -757778

$ws.Range("H80").Value = 3479.8
$ws.Range("I80").Value = 3649.5
$ws.Range("J80").Value = 3366.6667
$ws.Range("K80").Value = 3649.5
$ws.Range("L80").Value = 3366.6667
$ws.Range("M80").Value = -2651.5
$ws.Range("N80").Value = -5362.6667

$ws.Range("H83").Value = 3479.8
$ws.Range("I83").Value = 3649.5
$ws.Range("J83").Value = 3366.6667
$ws.Range("K83").Value = 18247.5
$ws.Range("L83").Value = 16833.3335
$ws.Range("M83").Value = -13255.5
$ws.Range("N83").Value = -26817.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3083.1667
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 1249.5
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 1249.5
$ws.Range("M22").Value = -3705
$ws.Range("N22").Value = -1839.5

$ws.Range("H27").Value = 3083.1667
$ws.Range("I27").Value = 4000
$ws.Range("J27").Value = 1249.5
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 1249.5
$ws.Range("M27").Value = -3893
$ws.Range("N27").Value = -1463.5

$ws.Range("H40").Value = 4599.1665
$ws.Range("J40").Value = 4774.25
$ws.Range("L40").Value = 4774.25
$ws.Range("N40").Value = -5046.25

$ws.Range("H41").Value = 19999
$ws.Range("I41").Value = 19999
$ws.Range("K41").Value = 19999
$ws.Range("M41").Value = -19561

$ws.Range("H46").Value = 2461.8333
$ws.Range("I46").Value = 1896.625
$ws.Range("J46").Value = 3592.25
$ws.Range("K46").Value = 1896.625
$ws.Range("L46").Value = 3592.25
$ws.Range("M46").Value = -1708.625
$ws.Range("N46").Value = -3968.25

$ws.Range("H82").Value = 1324.4
$ws.Range("I82").Value = 1100
$ws.Range("J82").Value = 1420.5714
$ws.Range("K82").Value = 1100
$ws.Range("L82").Value = 1420.5714
$ws.Range("M82").Value = -739
$ws.Range("N82").Value = -2142.5714

$ws.Range("H85").Value = 1324.4
$ws.Range("I85").Value = 1100
$ws.Range("J85").Value = 1420.5714
$ws.Range("K85").Value = 1100
$ws.Range("L85").Value = 1420.5714
$ws.Range("M85").Value = 148
$ws.Range("N85").Value = -3916.5714

$ws.Range("H93").Value = 1734.5555
$ws.Range("J93").Value = 1599.6666
$ws.Range("L93").Value = 1599.6666
$ws.Range("N93").Value = -4095.6666

$ws.Range("H100").Value = 3899.4
$ws.Range("J100").Value = 2799
$ws.Range("L100").Value = 2799
$ws.Range("N100").Value = -3881

$ws.Range("H132").Value = 3358.6
$ws.Range("I132").Value = 3169.5
$ws.Range("J132").Value = 3736.8
$ws.Range("K132").Value = 9508.5
$ws.Range("L132").Value = 11210.4
$ws.Range("M132").Value = -6978.5
$ws.Range("N132").Value = -16270.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1001095
$ws.Range("I81").Value = 1216.4445
$ws.Range("K81").Value = 2432.889
$ws.Range("M81").Value = -1371.889

$ws.Range("H84").Value = 1001095
$ws.Range("I84").Value = 1216.4445
$ws.Range("K84").Value = 12164.445
$ws.Range("M84").Value = -6860.445

$ws.Range("H96").Value = 1176
$ws.Range("I96").Value = 1018.3333
$ws.Range("K96").Value = 1018.3333
$ws.Range("M96").Value = 354.6667
